# Auto-generated edit script applying the Phoenix_Profits.xlsx cell updates
# described in the commit diff (scheduled runner refresh of market-price
# derived columns H:N across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 594852.4
$ws.Range("I69").Value = 631780.6
$ws.Range("J69").Value = 4000
$ws.Range("K69").Value = 1895341.8
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = -1894467.8
$ws.Range("N69").Value = -13748
$ws.Range("H72").Value = 594852.4
$ws.Range("I72").Value = 631780.6
$ws.Range("J72").Value = 4000
$ws.Range("K72").Value = 5686025.399999999
$ws.Range("L72").Value = 36000
$ws.Range("M72").Value = -5681657.399999999
$ws.Range("N72").Value = -44736
$ws.Range("H123").Value = 61249
$ws.Range("J123").Value = 61249
$ws.Range("L123").Value = 61249
$ws.Range("N123").Value = -71049
$ws.Range("H132").Value = 1994.3125
$ws.Range("I132").Value = 1528.5349
$ws.Range("K132").Value = 4585.6047
$ws.Range("M132").Value = -2055.6047
$ws.Range("H138").Value = 3451.328
$ws.Range("I138").Value = 1398
$ws.Range("K138").Value = 4194
$ws.Range("M138").Value = 946

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19601.137
$ws.Range("I32").Value = 18606.508
$ws.Range("K32").Value = 18606.508
$ws.Range("M32").Value = -18319.508
$ws.Range("H61").Value = 4785.186
$ws.Range("I61").Value = 3056.8147
$ws.Range("J61").Value = 7701.8125
$ws.Range("K61").Value = 3056.8147
$ws.Range("L61").Value = 7701.8125
$ws.Range("M61").Value = -2844.8147
$ws.Range("N61").Value = -8125.8125
$ws.Range("H74").Value = 2218.6875
$ws.Range("J74").Value = 3920.6
$ws.Range("L74").Value = 3920.6
$ws.Range("N74").Value = -5668.6
$ws.Range("H77").Value = 2218.6875
$ws.Range("J77").Value = 3920.6
$ws.Range("L77").Value = 19603
$ws.Range("N77").Value = -28339
$ws.Range("H122").Value = 78399.914
$ws.Range("I122").Value = 3079.5
$ws.Range("K122").Value = 9238.5
$ws.Range("M122").Value = -6788.5
$ws.Range("H132").Value = 10922.424
$ws.Range("I132").Value = 12400.556
$ws.Range("K132").Value = 37201.66800000001
$ws.Range("M132").Value = -34671.66800000001
$ws.Range("H136").Value = 4785.186
$ws.Range("I136").Value = 3056.8147
$ws.Range("J136").Value = 7701.8125
$ws.Range("K136").Value = 9170.444100000001
$ws.Range("L136").Value = 23105.4375
$ws.Range("M136").Value = -6620.444100000001
$ws.Range("N136").Value = -28205.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 91792.60000000001
$ws.Range("I86").Value = 1790.8572
$ws.Range("J86").Value = 301796.66
$ws.Range("K86").Value = 1790.8572
$ws.Range("L86").Value = 301796.66
$ws.Range("M86").Value = -667.8571999999999
$ws.Range("N86").Value = -304042.66
$ws.Range("H89").Value = 91792.60000000001
$ws.Range("I89").Value = 1790.8572
$ws.Range("J89").Value = 301796.66
$ws.Range("K89").Value = 8954.286
$ws.Range("L89").Value = 1508983.3
$ws.Range("M89").Value = -3338.286
$ws.Range("N89").Value = -1520215.3
$ws.Range("H128").Value = 999
$ws.Range("I128").Value = 999
$ws.Range("K128").Value = 2997
$ws.Range("M128").Value = -507
$ws.Range("H134").Value = 240796.67
$ws.Range("I134").Value = 1679.1852
$ws.Range("K134").Value = 5037.5556
$ws.Range("M134").Value = -2502.5556

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2638.4092
$ws.Range("I31").Value = 2376.9333
$ws.Range("K31").Value = 2376.9333
$ws.Range("M31").Value = -2081.9333
$ws.Range("H34").Value = 2638.4092
$ws.Range("I34").Value = 2376.9333
$ws.Range("K34").Value = 2376.9333
$ws.Range("M34").Value = -2174.9333
$ws.Range("H50").Value = 33100.945
$ws.Range("J50").Value = 35115.6
$ws.Range("L50").Value = 35115.6
$ws.Range("N50").Value = -36365.6
$ws.Range("H51").Value = 29541.666
$ws.Range("I51").Value = 28544
$ws.Range("J51").Value = 29826.715
$ws.Range("K51").Value = 28544
$ws.Range("L51").Value = 29826.715
$ws.Range("M51").Value = -27808
$ws.Range("N51").Value = -31298.715
$ws.Range("H58").Value = 1932.7059
$ws.Range("I58").Value = 1297.6154
$ws.Range("K58").Value = 1297.6154
$ws.Range("M58").Value = -1094.6154
$ws.Range("H59").Value = 45717.223
$ws.Range("J59").Value = 48931.875
$ws.Range("L59").Value = 48931.875
$ws.Range("N59").Value = -51221.875
$ws.Range("H61").Value = 29541.666
$ws.Range("I61").Value = 28544
$ws.Range("J61").Value = 29826.715
$ws.Range("K61").Value = 28544
$ws.Range("L61").Value = 29826.715
$ws.Range("M61").Value = -28196
$ws.Range("N61").Value = -30522.715
$ws.Range("H68").Value = 37555.707
$ws.Range("J68").Value = 38627.934
$ws.Range("L68").Value = 38627.934
$ws.Range("N68").Value = -40125.934
$ws.Range("H71").Value = 37555.707
$ws.Range("J71").Value = 38627.934
$ws.Range("L71").Value = 115883.802
$ws.Range("N71").Value = -123371.802
$ws.Range("H94").Value = 1304.5834
$ws.Range("J94").Value = 663.1429000000001
$ws.Range("L94").Value = 663.1429000000001
$ws.Range("N94").Value = -1565.1429
$ws.Range("H132").Value = 386280.72
$ws.Range("I132").Value = 434117.66
$ws.Range("J132").Value = 27503.5
$ws.Range("K132").Value = 1302352.98
$ws.Range("L132").Value = 82510.5
$ws.Range("M132").Value = -1299822.98
$ws.Range("N132").Value = -87570.5
$ws.Range("H134").Value = 2554.4666
$ws.Range("I134").Value = 1770.591
$ws.Range("J134").Value = 4710.125
$ws.Range("K134").Value = 5311.772999999999
$ws.Range("L134").Value = 14130.375
$ws.Range("M134").Value = -2776.772999999999
$ws.Range("N134").Value = -19200.375
$ws.Range("H136").Value = 1932.7059
$ws.Range("I136").Value = 1297.6154
$ws.Range("K136").Value = 3892.8462
$ws.Range("M136").Value = -1342.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2899.7144
$ws.Range("J131").Value = 3704.5
$ws.Range("L131").Value = 11113.5
$ws.Range("N131").Value = -21193.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""
$ws.Range("H80").Value = 8043.8667
$ws.Range("J80").Value = 11077.223
$ws.Range("L80").Value = 11077.223
$ws.Range("N80").Value = -13073.223
$ws.Range("H83").Value = 8043.8667
$ws.Range("J83").Value = 11077.223
$ws.Range("L83").Value = 55386.115
$ws.Range("N83").Value = -65370.115
$ws.Range("H122").Value = 4560
$ws.Range("I122").Value = 8500
$ws.Range("K122").Value = 25500
$ws.Range("M122").Value = -23050

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2415.15
$ws.Range("I61").Value = 2189.4856
$ws.Range("K61").Value = 2189.4856
$ws.Range("M61").Value = -1987.4856
$ws.Range("H113").Value = 2415.15
$ws.Range("I113").Value = 2189.4856
$ws.Range("K113").Value = 2189.4856
$ws.Range("M113").Value = -19.48559999999998
$ws.Range("H122").Value = 3234.0322
$ws.Range("I122").Value = 3196.7585
$ws.Range("J122").Value = 3774.5
$ws.Range("K122").Value = 9590.2755
$ws.Range("L122").Value = 11323.5
$ws.Range("M122").Value = -7140.2755
$ws.Range("N122").Value = -16223.5
$ws.Range("H132").Value = 2987.9783
$ws.Range("I132").Value = 2661.9
$ws.Range("J132").Value = 5161.8335
$ws.Range("K132").Value = 7985.700000000001
$ws.Range("L132").Value = 15485.5005
$ws.Range("M132").Value = -5455.700000000001
$ws.Range("N132").Value = -20545.5005
$ws.Range("H136").Value = 3310.8604
$ws.Range("I136").Value = 2885.4333
$ws.Range("K136").Value = 8656.2999
$ws.Range("M136").Value = -6106.2999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""
$ws.Range("H81").Value = 2830.35
$ws.Range("I81").Value = 952.9167
$ws.Range("J81").Value = 5646.5
$ws.Range("K81").Value = 1905.8334
$ws.Range("L81").Value = 11293
$ws.Range("M81").Value = -844.8334
$ws.Range("N81").Value = -13415
$ws.Range("H84").Value = 2830.35
$ws.Range("I84").Value = 952.9167
$ws.Range("J84").Value = 5646.5
$ws.Range("K84").Value = 9529.166999999999
$ws.Range("L84").Value = 56465
$ws.Range("M84").Value = -4225.166999999999
$ws.Range("N84").Value = -67073
$ws.Range("H122").Value = 4970.1763
$ws.Range("J122").Value = 5028
$ws.Range("L122").Value = 15084
$ws.Range("N122").Value = -19984
$ws.Range("H126").Value = 31581.857
$ws.Range("I126").Value = 36096.75
$ws.Range("K126").Value = 108290.25
$ws.Range("M126").Value = -105820.25
$ws.Range("H132").Value = 3476.42
$ws.Range("I132").Value = 2340.1724
$ws.Range("J132").Value = 5045.524
$ws.Range("K132").Value = 7020.5172
$ws.Range("L132").Value = 15136.572
$ws.Range("M132").Value = -4490.5172
$ws.Range("N132").Value = -20196.572
$ws.Range("H136").Value = 2090.1892
$ws.Range("I136").Value = 1698.3438
$ws.Range("J136").Value = 4598
$ws.Range("K136").Value = 5095.0314
$ws.Range("L136").Value = 13794
$ws.Range("M136").Value = -2545.0314
$ws.Range("N136").Value = -18894
